# Add "team record" columns (Wins/Losses/Ties) to the MIA_2018 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, matching the style of the existing header row (A1:AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Every player row gets the same season team record.
$ws.Range("AD2:AD53").Value = 63
$ws.Range("AE2:AE53").Value = 98
$ws.Range("AF2:AF53").Value = 0
